# Apply odds updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("F3").Value = 2.12
$ws.Range("G3").Value = 2.52

# Row 4
$ws.Range("F4").Value = 1.66
$ws.Range("I4").Value = 5.2
$ws.Range("J4").Value = 4.3
$ws.Range("K4").Value = 5.6

# Row 5
$ws.Range("G5").Value = 1.52
$ws.Range("H5").Value = 8.6
$ws.Range("J5").Value = 4.4
$ws.Range("P5").Value = 2.16
$ws.Range("Q5").Value = 1.68

# Row 6
$ws.Range("H6").Value = 2.52
$ws.Range("J6").Value = 3.4
$ws.Range("P6").Value = 2.02
$ws.Range("Q6").Value = 1.76

# Row 7
$ws.Range("F7").Value = 1.74
$ws.Range("I7").Value = 5.3
$ws.Range("J7").Value = 3.2
$ws.Range("K7").Value = 85
$ws.Range("Q7").Value = 1.71

# Row 8
$ws.Range("G8").Value = 1.71
$ws.Range("K8").Value = 5.1
$ws.Range("P8").Value = 2.4
$ws.Range("Q8").Value = 1.49

# Row 9
$ws.Range("Q9").Value = 2.04
